$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column F (modBusIndex)
$ws.Columns("F").Delete()

# Update autoSub sample value in row 3 (now column R after the deletion) from "0/1" to "1"
$ws.Range("R3").Value = "1"
